$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 28.125
$ws.Range("I2").Value = 27.826086
$ws.Range("J2").Value = 35
$ws.Range("K2").Value = 27.826086
$ws.Range("L2").Value = 35
$ws.Range("M2").Value = 85.173914
$ws.Range("N2").Value = -261

# Row 18
$ws.Range("H18").Value = 125011000
$ws.Range("I18").Value = 200001400
$ws.Range("J18").Value = 26999.334
$ws.Range("K18").Value = 200001400
$ws.Range("L18").Value = 26999.334
$ws.Range("M18").Value = -200001116
$ws.Range("N18").Value = -27567.334

# Row 19
$ws.Range("H19").Value = 4690.5454
$ws.Range("I19").Value = 4562.6665
$ws.Range("K19").Value = 4562.6665
$ws.Range("M19").Value = -4387.6665

# Row 29
$ws.Range("H29").Value = 1499.5
$ws.Range("I29").Value = 999
$ws.Range("K29").Value = 2997
$ws.Range("M29").Value = -2716

# Row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

# Row 43
$ws.Range("H43").Value = 689992
$ws.Range("J43").Value = 1032500.5
$ws.Range("L43").Value = 1032500.5
$ws.Range("N43").Value = -1032638.5

# Row 74
$ws.Range("H74").Value = 48393120
$ws.Range("I74").Value = 75002710
$ws.Range("J74").Value = 12045.728
$ws.Range("K74").Value = 75002710
$ws.Range("L74").Value = 12045.728
$ws.Range("M74").Value = -75001774
$ws.Range("N74").Value = -13917.728

# Row 77
$ws.Range("H77").Value = 48393120
$ws.Range("I77").Value = 75002710
$ws.Range("J77").Value = 12045.728
$ws.Range("K77").Value = 375013550
$ws.Range("L77").Value = 60228.64
$ws.Range("M77").Value = -375008870
$ws.Range("N77").Value = -69588.64

# Row 113
$ws.Range("H113").Value = 67141110
$ws.Range("J113").Value = 93765544
$ws.Range("L113").Value = 93765544
$ws.Range("N113").Value = -93772052

# Row 132
$ws.Range("H132").Value = 1514.4667
$ws.Range("I132").Value = 1479.3572
$ws.Range("K132").Value = 4438.071599999999
$ws.Range("M132").Value = -1908.071599999999

# Row 138
$ws.Range("H138").Value = 1591435.9
$ws.Range("I138").Value = 2051.476
$ws.Range("J138").Value = 2386128
$ws.Range("K138").Value = 6154.428
$ws.Range("L138").Value = 7158384
$ws.Range("M138").Value = -1014.428
$ws.Range("N138").Value = -7168664


$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 76925470
$ws.Range("I2").Value = 2045.4445
$ws.Range("K2").Value = 2045.4445
$ws.Range("M2").Value = -1932.4445

# Row 32
$ws.Range("H32").Value = 4173953.8
$ws.Range("I32").Value = 4550904.5
$ws.Range("K32").Value = 4550904.5
$ws.Range("M32").Value = -4550617.5

# Row 45
$ws.Range("H45").Value = 3085.5
$ws.Range("I45").Value = 3300
$ws.Range("J45").Value = 3014
$ws.Range("K45").Value = 3300
$ws.Range("L45").Value = 3014
$ws.Range("M45").Value = -2923
$ws.Range("N45").Value = -3768

# Row 61
$ws.Range("H61").Value = 27783078
$ws.Range("I61").Value = 2206.4211
$ws.Range("J61").Value = 58832290
$ws.Range("K61").Value = 2206.4211
$ws.Range("L61").Value = 58832290
$ws.Range("M61").Value = -1994.4211
$ws.Range("N61").Value = -58832714

# Row 74
$ws.Range("H74").Value = 24692.2
$ws.Range("I74").Value = 31730.94
$ws.Range("K74").Value = 31730.94
$ws.Range("M74").Value = -30856.94

# Row 77
$ws.Range("H77").Value = 24692.2
$ws.Range("I77").Value = 31730.94
$ws.Range("K77").Value = 158654.7
$ws.Range("M77").Value = -154286.7

# Row 116
$ws.Range("H116").Value = 76925470
$ws.Range("I116").Value = 2045.4445
$ws.Range("K116").Value = 2045.4445
$ws.Range("M116").Value = 248.5554999999999

# Row 132
$ws.Range("H132").Value = 7230.9756
$ws.Range("I132").Value = 5424.0386
$ws.Range("J132").Value = 10363
$ws.Range("K132").Value = 16272.1158
$ws.Range("L132").Value = 31089
$ws.Range("M132").Value = -13742.1158
$ws.Range("N132").Value = -36149

# Row 136
$ws.Range("H136").Value = 27783078
$ws.Range("I136").Value = 2206.4211
$ws.Range("J136").Value = 58832290
$ws.Range("K136").Value = 6619.263300000001
$ws.Range("L136").Value = 176496870
$ws.Range("M136").Value = -4069.263300000001
$ws.Range("N136").Value = -176501970


$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 76925470
$ws.Range("I3").Value = 2045.4445
$ws.Range("K3").Value = 2045.4445
$ws.Range("M3").Value = -1931.4445

# Row 22
$ws.Range("H22").Value = 13889113
$ws.Range("I22").Value = 13889113
$ws.Range("K22").Value = 13889113
$ws.Range("M22").Value = -13888940

# Row 94
$ws.Range("H94").Value = 1962.5883
$ws.Range("I94").Value = 1390.6
$ws.Range("J94").Value = 6252.5
$ws.Range("K94").Value = 1390.6
$ws.Range("L94").Value = 6252.5
$ws.Range("M94").Value = -939.5999999999999
$ws.Range("N94").Value = -7154.5

# Row 99
$ws.Range("H99").Value = 4786499
$ws.Range("I99").Value = 1493.7858
$ws.Range("K99").Value = 1493.7858
$ws.Range("M99").Value = 4.214199999999892

# Row 105
$ws.Range("H105").Value = 1669.9375
$ws.Range("I105").Value = 1512.909
$ws.Range("J105").Value = 2015.4
$ws.Range("K105").Value = 1512.909
$ws.Range("L105").Value = 2015.4
$ws.Range("M105").Value = 234.0909999999999
$ws.Range("N105").Value = -5509.4

# Row 107
$ws.Range("H107").Value = 53627060
$ws.Range("I107").Value = 75073890
$ws.Range("J107").Value = 9998.833000000001
$ws.Range("K107").Value = 75073890
$ws.Range("L107").Value = 9998.833000000001
$ws.Range("M107").Value = -75071970
$ws.Range("N107").Value = -13838.833


$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 772.6667
$ws.Range("I22").Value = 766.6667
$ws.Range("J22").Value = 778.6667
$ws.Range("K22").Value = 766.6667
$ws.Range("L22").Value = 778.6667
$ws.Range("M22").Value = -416.6667
$ws.Range("N22").Value = -1478.6667

# Row 51
$ws.Range("H51").Value = 45316.668
$ws.Range("J51").Value = 45316.668
$ws.Range("L51").Value = 45316.668
$ws.Range("N51").Value = -46788.668

# Row 61
$ws.Range("H61").Value = 45316.668
$ws.Range("J61").Value = 45316.668
$ws.Range("L61").Value = 45316.668
$ws.Range("N61").Value = -46012.668

# Row 94
$ws.Range("H94").Value = 1304.6471
$ws.Range("I94").Value = 1574
$ws.Range("J94").Value = 1157.7273
$ws.Range("K94").Value = 1574
$ws.Range("L94").Value = 1157.7273
$ws.Range("M94").Value = -1123
$ws.Range("N94").Value = -2059.7273

# Row 107
$ws.Range("H107").Value = 2190.3044
$ws.Range("I107").Value = 664.8333
$ws.Range("J107").Value = 2728.7058
$ws.Range("K107").Value = 664.8333
$ws.Range("L107").Value = 2728.7058
$ws.Range("M107").Value = 1255.1667
$ws.Range("N107").Value = -6568.7058

# Row 141
$ws.Range("H141").Value = 378184.34
$ws.Range("J141").Value = 378184.34
$ws.Range("L141").Value = 378184.34
$ws.Range("N141").Value = -388544.34


$ws = $wb.Worksheets.Item("CUL")
# Row 47
$ws.Range("H47").Value = 700
$ws.Range("I47").Value = 700
$ws.Range("K47").Value = 2100
$ws.Range("M47").Value = -1669

# Row 92
$ws.Range("H92").Value = 12821937
$ws.Range("I92").Value = 774.5
$ws.Range("K92").Value = 2323.5
$ws.Range("M92").Value = -1075.5

# Row 107
$ws.Range("H107").Value = 28572310
$ws.Range("I107").Value = 700
$ws.Range("K107").Value = 2100
$ws.Range("M107").Value = -180

# Row 139
$ws.Range("H139").Value = 59109.445
$ws.Range("I139").Value = 65248.125
$ws.Range("K139").Value = 195744.375
$ws.Range("M139").Value = -190604.375


$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 406308.06
$ws.Range("I70").Value = 731567.5600000001
$ws.Range("K70").Value = 731567.5600000001
$ws.Range("M70").Value = -731297.5600000001

# Row 73
$ws.Range("H73").Value = 406308.06
$ws.Range("I73").Value = 731567.5600000001
$ws.Range("K73").Value = 731567.5600000001
$ws.Range("M73").Value = -730631.5600000001

# Row 80
$ws.Range("H80").Value = 5093.9287
$ws.Range("I80").Value = 3005.5715
$ws.Range("K80").Value = 3005.5715
$ws.Range("M80").Value = -2007.5715

# Row 83
$ws.Range("H83").Value = 5093.9287
$ws.Range("I83").Value = 3005.5715
$ws.Range("K83").Value = 15027.8575
$ws.Range("M83").Value = -10035.8575


$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 10401695
$ws.Range("I46").Value = 6897474
$ws.Range("J46").Value = 12348484
$ws.Range("K46").Value = 6897474
$ws.Range("L46").Value = 12348484
$ws.Range("M46").Value = -6897286
$ws.Range("N46").Value = -12348860

# Row 68
$ws.Range("H68").Value = 7397
$ws.Range("J68").Value = 7397
$ws.Range("L68").Value = 7397
$ws.Range("N68").Value = -8895

# Row 71
$ws.Range("H71").Value = 7397
$ws.Range("J71").Value = 7397
$ws.Range("L71").Value = 36985
$ws.Range("N71").Value = -44473

# Row 93
$ws.Range("H93").Value = 778.8889
$ws.Range("I93").Value = 790.5
$ws.Range("K93").Value = 790.5
$ws.Range("M93").Value = 457.5


$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 2624.5
$ws.Range("I96").Value = 2998
$ws.Range("J96").Value = 2500
$ws.Range("K96").Value = 2998
$ws.Range("L96").Value = 2500
$ws.Range("M96").Value = -1625
$ws.Range("N96").Value = -5246

# Row 100
$ws.Range("H100").Value = 612.6818
$ws.Range("I100").Value = 383.45456
$ws.Range("J100").Value = 841.9091
$ws.Range("K100").Value = 766.90912
$ws.Range("L100").Value = 1683.8182
$ws.Range("M100").Value = -225.90912
$ws.Range("N100").Value = -2765.8182

# Row 122
$ws.Range("H122").Value = 113361.445
$ws.Range("I122").Value = 139557.77
$ws.Range("J122").Value = 4833.857
$ws.Range("K122").Value = 418673.3099999999
$ws.Range("L122").Value = 14501.571
$ws.Range("M122").Value = -416223.3099999999
$ws.Range("N122").Value = -19401.571

# Row 132
$ws.Range("H132").Value = 10688.167
$ws.Range("I132").Value = 18509.834
$ws.Range("J132").Value = 2866.5
$ws.Range("K132").Value = 55529.50199999999
$ws.Range("L132").Value = 8599.5
$ws.Range("M132").Value = -52999.50199999999
$ws.Range("N132").Value = -13659.5

# Row 136
$ws.Range("H136").Value = 24051238
$ws.Range("I136").Value = 58824684
$ws.Range("J136").Value = 405296
$ws.Range("K136").Value = 176474052
$ws.Range("L136").Value = 1215888
$ws.Range("M136").Value = -176471502
$ws.Range("N136").Value = -1220988

